# "2 commit del dia"
# Fill in the "Entrega" (delivery) denomination counts for the first
# cash-count block (columns M:Q, dates 43254/43261/43268) in column P,
# rows 4-13. Column Q already holds =SUM(M*P) formulas that recalculate
# automatically once P is populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P4").Value = 5
$ws.Range("P5").Value = 19
$ws.Range("P6").Value = 15
$ws.Range("P7").Value = 29
$ws.Range("P8").Value = 13
$ws.Range("P9").Value = 9
$ws.Range("P10").Value = 4
$ws.Range("P11").Value = 0
$ws.Range("P12").Value = 2
$ws.Range("P13").Value = 0

$ws.Range("P13").Select()
